$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Version and Date, insert Jurisdiction row ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row before row 11 (currently "Description") to hold "Jurisdiction"
$meta.Rows.Item(11).Insert()
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# --- Sheet "Elements": add constraint text to Authorization.typeId row ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AJ5").Value = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}
"
